$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 47 (shifts existing rows 47-76 down to 48-77,
# carrying along formatting such as the date style on column D).
$ws.Rows.Item(47).Insert()

# Column A, B, C, E, F, G, R are identical for every record in this block
# (same market/region/category). Set the same literal values in the
# freshly inserted row 47.
$ws.Cells.Item(47,1).Value  = 11                            # A: Mercado ID
$ws.Cells.Item(47,2).Value  = "Vega Monumental Concepción"  # B: Mercado
$ws.Cells.Item(47,3).Value  = "Bíobío"                       # C: Región
$ws.Cells.Item(47,5).Value  = 8                              # E: Codreg
$ws.Cells.Item(47,6).Value  = 100112024                      # F: Categoría ID
$ws.Cells.Item(47,7).Value  = "Choclo"                        # G: Categoría
$ws.Cells.Item(47,18).Value = "Hortaliza"                     # R: Clasificación

# New record data for the inserted row.
$ws.Cells.Item(47,4).Value  = 44603                    # D: Fecha
$ws.Cells.Item(47,8).Value  = "Choclero"               # H: Variedad
$ws.Cells.Item(47,9).Value  = "Primera"                # I: Calidad
$ws.Cells.Item(47,10).Value = 3000                     # J: Volumen
$ws.Cells.Item(47,11).Value = 90                       # K: Precio mínimo
$ws.Cells.Item(47,12).Value = 150                      # L: Precio máximo
$ws.Cells.Item(47,13).Value = 120                      # M: Precio promedio ponderado
$ws.Cells.Item(47,14).Value = "$/unidad"               # N: Unidad de comercialización
$ws.Cells.Item(47,15).Value = "Región Metropolitana"   # O: Origen
$ws.Cells.Item(47,16).Value = 120                      # P: Precio $/Kg
$ws.Cells.Item(47,17).Value = 1                        # Q: Kg o Unidades
